$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Obtido" column (E): "não" -> "sim" for these parts
$ws.Range("E3").Value  = "sim"
$ws.Range("E6").Value  = "sim"
$ws.Range("E7").Value  = "sim"
$ws.Range("E13").Value = "sim"
$ws.Range("E23").Value = "sim"
$ws.Range("E24").Value = "sim"

# "Obtenção" column (F): update progress status
$ws.Range("F14").Value = "incompleto"
$ws.Range("F16").Value = "incompleto"
$ws.Range("F29").Value = "concluído"
$ws.Range("F31").Value = "incompleto"

# Update the active selection / scroll position to match the latest edit location
$ws.Range("F16").Select()
